# Weekly data refresh: a new record (most recent survey date) is added at
# the top of the price table (row 29, right after the 27 rows belonging to
# other markets/varieties that precede this "Arveja Verde" subset), pushing
# every existing "Arveja Verde" record down by one row. The oldest record
# (previously the last row, 112) ends up in the new last row, 113.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 29; this shifts rows 29:112 down to
# rows 30:113 (dimension grows from A1:R112 to A1:R113 automatically).
$ws.Rows(29).Insert()

# Populate the newly inserted row 29 with this week's new record.
$ws.Range("A29").Value2 = 9
$ws.Range("B29").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C29").Value2 = "Metropolitana"
$ws.Range("D29").Value2 = 44659
$ws.Range("E29").Value2 = 13
$ws.Range("F29").Value2 = 100112022
$ws.Range("G29").Value2 = "Arveja Verde"
$ws.Range("H29").Value2 = "Sin especificar"
$ws.Range("I29").Value2 = "Primera"
$ws.Range("J29").Value2 = 52
$ws.Range("K29").Value2 = 23000
$ws.Range("L29").Value2 = 25000
$ws.Range("M29").Value2 = 24000
$ws.Range("N29").Value2 = '$/saco 25 kilos'
$ws.Range("O29").Value2 = "Carahue"
$ws.Range("P29").Value2 = 960
$ws.Range("Q29").Value2 = 25
$ws.Range("R29").Value2 = "Hortaliza"

# Make sure the D column of the new row keeps the date number format used
# by the rest of the column (style carried over from the row above on
# insert, but set it explicitly to be safe).
$ws.Range("D29").NumberFormat = $ws.Range("D30").NumberFormat
